$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.230.29'
$ws.Range('E2').Value = '  +3.57%  '
$ws.Range('D3').Value = '3.060.00'
$ws.Range('E3').Value = '  +3.41%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '521.12'
$ws.Range('E5').Value = '  +4.91%  '
$ws.Range('D6').Value = '141.64'
$ws.Range('E6').Value = '  +6.66%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  +5.13%  '
$ws.Range('D9').Value = '7.52'
$ws.Range('E9').Value = '  +4.00%  '
$ws.Range('D10').Value = '0.112'
$ws.Range('E10').Value = '  +6.39%  '
$ws.Range('D11').Value = '0.369'
$ws.Range('E11').Value = '  +6.09%  '
$ws.Range('D12').Value = '3.504.32'
$ws.Range('E12').Value = '  +1.08%  '
$ws.Range('E13').Value = '  +2.83%  '
$ws.Range('D14').Value = '26.84'
$ws.Range('E14').Value = '  +7.30%  '
$ws.Range('D15').Value = '0.0000171'
$ws.Range('E15').Value = '  +15.73%  '
$ws.Range('D16').Value = '58.165.18'
$ws.Range('E16').Value = '  +3.38%  '
$ws.Range('D17').Value = '6.26'
$ws.Range('E17').Value = '  +10.62%  '
$ws.Range('D18').Value = '3.051.39'
$ws.Range('E18').Value = '  +2.95%  '
$ws.Range('D19').Value = '13.08'
$ws.Range('E19').Value = '  +7.15%  '
$ws.Range('D20').Value = '8.15'
$ws.Range('E20').Value = '  +6.12%  '
$ws.Range('D21').Value = '337.78'
$ws.Range('E21').Value = '  +4.51%  '
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D22').Value = '5.78'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '0.502'
$ws.Range('E24').Value = '  +7.64%  '
$ws.Range('D25').Value = '65.20'
$ws.Range('E25').Value = '  +5.77%  '
$ws.Range('D26').Value = '0.169'
$ws.Range('E26').Value = '  +4.23%  '
$ws.Range('D27').Value = '0.0₃0961'
$ws.Range('E27').Value = '  +9.41%  '
$ws.Range('D28').Value = '0.993'
$ws.Range('E28').Value = '  -0.54%  '
$ws.Range('D29').Value = '6.91'
$ws.Range('E29').Value = '  +8.39%  '
$ws.Range('D30').Value = '7.55'
$ws.Range('E30').Value = '  +12.35%  '
$ws.Range('D31').Value = '1.84'
$ws.Range('E31').Value = '  +6.60%  '
$ws.Range('D32').Value = '1.22'
$ws.Range('E32').Value = '  +5.96%  '
$ws.Range('D33').Value = '21.14'
$ws.Range('E33').Value = '  +5.25%  '
$ws.Range('D34').Value = '157.70'
$ws.Range('E34').Value = '  +1.59%  '
$ws.Range('D35').Value = '4.78'
$ws.Range('E35').Value = '  +9.04%  '
$ws.Range('D36').Value = '5.96'
$ws.Range('E36').Value = '  +8.54%  '
$ws.Range('E37').Value = '  +3.41%  '
$ws.Range('D38').Value = '25.51'
$ws.Range('E38').Value = '  +13.34%  '
$ws.Range('D39').Value = '0.0695'
$ws.Range('E39').Value = '  +4.68%  '
$ws.Range('D40').Value = '3.093.63'
$ws.Range('E40').Value = '  +3.38%  '
$ws.Range('D41').Value = '37.77'
$ws.Range('E41').Value = '  +5.87%  '
$ws.Range('D42').Value = '3.91'
$ws.Range('E42').Value = '  +11.15%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.667'
$ws.Range('E43').Value = '  +4.89%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('D45').Value = '2.335.08'
$ws.Range('E45').Value = '  +5.36%  '
$ws.Range('D46').Value = '1.46'
$ws.Range('E46').Value = '  +5.98%  '
$ws.Range('E47').Value = '  +4.92%  '
$ws.Range('D48').Value = '6.07'
$ws.Range('E48').Value = '  +6.29%  '
$ws.Range('E49').Value = '  +4.02%  '
$ws.Range('D50').Value = '19.84'
$ws.Range('E50').Value = '  +6.37%  '
$ws.Range('D51').Value = '1.87'
$ws.Range('E51').Value = '  -0.81%  '
